$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(121, 2281, 101, 1225, 51, 3658, 51, 50, 101, 24459, 16430, 4371, 0, 20088, 0, 0, 1326, 0),
    @(122, 2281, 101, 1225, 51, 3658, 51, 50, 101, 24459, 16430, 4371, 0, 20088, 0, 0, 1326, 0)
)

$startRow = 122
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
    # Column S (19) stays empty, but still materialize the cell node,
    # matching the blank-but-present "S" cells used throughout the sheet.
    $ws.Cells.Item($row, 19).Style = "Normal"
}
